$d = $word.ActiveDocument

# --- Step 1: remove the "License Information" Heading2 paragraph (currently paragraph 4) ---
$licPara = $d.Paragraphs.Item(4)
if ($licPara.Range.Text.TrimEnd([char]13) -eq "License Information") {
    $licPara.Range.Delete()
}

# --- Step 2: remove the "This PDF version is provided under the same license." paragraph
#             (currently paragraph 5, right after the big license paragraph) ---
$pdfPara = $d.Paragraphs.Item(5)
if ($pdfPara.Range.Text.TrimEnd([char]13) -eq "This PDF version is provided under the same license.") {
    $pdfPara.Range.Delete()
}

# --- Step 3: rewrite the big "Translation Questions (unfoldingWord) is based on..." paragraph
#             (now paragraph 4) with the new resource-data text ---
$p = $d.Paragraphs.Item(4)
$r = $p.Range
# Trim the trailing paragraph mark off the range so we only clear the paragraph's content.
$r.End = $r.End - 1
$r.Text = ""
$paraStart = $r.Start

$boldText = "unfoldingWord® Translation Questions"
$rest = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. " + `
        "unfoldingWord® Translation Questions" + `
        " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from " + `
        "unfoldingWord® Translation Questions" + `
        " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

$r.InsertAfter($boldText)
$r.Collapse(0)
$r.InsertAfter($rest)
$r.Collapse(0)

# Bold just the leading "unfoldingWord® Translation Questions" run.
$boldRange = $d.Range($paraStart, $paraStart + $boldText.Length)
$boldRange.Font.Bold = $true

Write-Host "Final paragraph text:"
Write-Host $d.Paragraphs.Item(4).Range.Text
